# Nowcasts 2025Q4 refresh: roll the rolling date window forward by one
# quarter (rows 2-7 now cover 2025-09-30 .. 2025-12-15 instead of
# 2025-06-30 .. 2025-09-15) and update the Prognose/Revision figures that
# go with the new window. A few column widths were also nudged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width adjustments -------------------------------------------------
$ws.Columns.Item(4).ColumnWidth  = 15.64453125
$ws.Columns.Item(6).ColumnWidth  = 15.77734375
$ws.Columns.Item(7).ColumnWidth  = 15.24609375
$ws.Columns.Item(8).ColumnWidth  = 15.24609375
$ws.Columns.Item(10).ColumnWidth = 15.046875

# --- row labels (column A): text dates, not real Excel date values ----------
# Use a leading apostrophe so Excel stores these as text (matching the
# original workbook, where column A is a shared string, not a date serial),
# then reset the style so no stray number-format style sticks to the cell.
function Set-TextCell($cell, $text) {
    $r = $ws.Range($cell)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

Set-TextCell "A2" "2025-09-30"
Set-TextCell "A3" "2025-10-15"
Set-TextCell "A4" "2025-10-30"
Set-TextCell "A5" "2025-11-15"
Set-TextCell "A6" "2025-11-30"
Set-TextCell "A7" "2025-12-15"

# --- refreshed Prognose / Revision figures -----------------------------------
$ws.Range("B2").Value = 0.084495522222524294
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("B3").Value = 0.24558823178915448
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.1232315930075218
$ws.Range("E3").Value = 0.00085603867194212579
$ws.Range("F3").Value = 0.0067299272148173273
$ws.Range("G3").Value = 0.013826063441293054
$ws.Range("H3").Value = -0.0022484857818651197
$ws.Range("I3").Value = -0.00050086909632624118
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.019198442109247232

$ws.Range("B4").Value = 0.53112229876576733
$ws.Range("C4").Value = 0.14126037339871611
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.00011359290670087254
$ws.Range("F4").Value = -0.0000031943819736848092
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.0017237540287851098
$ws.Range("I4").Value = 0.085955787734916911
$ws.Range("J4").Value = 0.06340699436369357
$ws.Range("K4").Value = -0.0032485472032540985

$ws.Range("B5").Value = 0.50706722219268607
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.00016365178291100416
$ws.Range("E5").Value = 0.013489427126714531
$ws.Range("F5").Value = -0.019908257800860866
$ws.Range("G5").Value = -0.0048478338308531443
$ws.Range("H5").Value = -0.0023081088075249433
$ws.Range("I5").Value = 0.0045201379975708068
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.015164093041038607

$ws.Range("B6").Value = 0.38538667896126755
$ws.Range("C6").Value = -0.067180841035158062
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.000068071724188183792
$ws.Range("F6").Value = -0.0011364902590309705
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.0075564526926717229
$ws.Range("I6").Value = -0.04603193863014432
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.00015710766139831467

$ws.Range("B7").Value = 0.4499619135115851
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.11917181044606817
$ws.Range("E7").Value = -0.025250252677014837
$ws.Range("F7").Value = -0.004034887095550806
$ws.Range("G7").Value = -0.0017518946189426642
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.023559541504242332
